# "aggiunta una colonna nella tabella"
# Insert a new calculated column "Totale in binario" into table Tabella4,
# positioned immediately before the existing last column "Totale in decimale"
# (which shifts one column to the right, from I to J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

$lastColIndex = $tbl.ListColumns.Count          # 8 -> "Totale in decimale"
$oldLastCol = $tbl.ListColumns.Item($lastColIndex)
$oldHeaderRange = $oldLastCol.Range              # header + data, currently column I
$oldDataRange = $oldLastCol.DataBodyRange        # I2:I17

$firstDataRow = $oldDataRange.Row
$lastDataRow = $firstDataRow + $oldDataRange.Rows.Count - 1
$headerRow = $oldHeaderRange.Row
$oldColNumber = $oldHeaderRange.Column            # 9 == column I

# Capture old formulas + cached values before anything is moved/overwritten.
$decFormula = "=(SUM(SUM(PRODUCT(Tabella4[[#This Row],[A1]]*2) + (PRODUCT(Tabella4[[#This Row],[A0]] * 1))) + SUM(PRODUCT(Tabella4[[#This Row],[B1]]*2) + (PRODUCT(Tabella4[[#This Row],[B0]] * 1)))))"
$binFormula = "=DEC2BIN(SUM(SUM(PRODUCT(Tabella4[[#This Row],[A1]]*2) + (PRODUCT(Tabella4[[#This Row],[A0]] * 1))) + SUM(PRODUCT(Tabella4[[#This Row],[B1]]*2) + (PRODUCT(Tabella4[[#This Row],[B0]] * 1)))))"

# 1) Append a brand-new table column (the engine always appends at the end,
#    so it lands one to the right of "Totale in decimale" -> becomes column J).
$newCol = $tbl.ListColumns.Add()
$newColNumber = $newCol.Range.Column               # 10 == column J

# 2) Re-label the headers so that J keeps the name "Totale in decimale" and
#    I becomes the new "Totale in binario" column (rename I first so there is
#    never a moment with two columns sharing the same name).
$ws.Cells.Item($headerRow, $oldColNumber).Value = "Totale in binario"
$ws.Cells.Item($headerRow, $newColNumber).Value = "Totale in decimale"

# 3) Move the old decimal-total formulas/values into the new J column.
$newDataRange = $ws.Range($ws.Cells.Item($firstDataRow, $newColNumber), $ws.Cells.Item($lastDataRow, $newColNumber))
$newDataRange.Formula = $decFormula
$newDataRange.NumberFormat = "0"
$newDataRange.NumberFormat = "General"

# 4) Put the new DEC2BIN formula into column I (old decimal column slot).
$oldDataRange.Formula = $binFormula

# 5) Formatting: the new "Totale in binario" column data is right-aligned.
$oldDataRange.HorizontalAlignment = -4152   # xlRight
$oldDataRange.NumberFormat = "General"

# 6) Column widths to roughly match the new layout.
$ws.Columns.Item($oldColNumber).ColumnWidth = 17.21875
$ws.Columns.Item($newColNumber).ColumnWidth = 19.44140625

# 7) Leave the selection where the edit session ended.
$ws.Range("H20").Select()
